# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gains a new blank column at N (pushing the
# existing "Late" / "heading" / "Outstanding" columns one slot to the right,
# N->O, O->P, P->Q) and becomes the active/selected sheet (with R8 selected),
# matching the recorded diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N; existing N/O/P data (and the
# column formatting that goes with it) shifts right to O/P/Q.
$ws.Columns("N").Insert()

# The freshly inserted column picks up a plain custom width (no bestFit),
# matching column M's rendered width of 11 character-units.
$ws.Columns("N").ColumnWidth = 10.17

# Make "Repayment schedule" the active sheet/tab, with R8 selected.
$ws.Activate() | Out-Null
$ws.Range("R8").Select() | Out-Null
